$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column S: year 2022 header, copying format from R4 (style s="7")
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# Add new column S: data value for row 5, copying format from R5 (style s="15")
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 3.4

# Update existing values in row 5 (P, Q, R)
$ws.Range("P5").Value = 4.4000000000000004
$ws.Range("Q5").Value = 2.9
$ws.Range("R5").Value = 3.2

# Update the selected cell in the sheet view
$ws.Range("T4").Select()

